# Update OpenAI news (with translation)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of the data (rows 2 and 3), pushing the
# existing data rows down by two (old row 2 -> new row 4, etc.)
$ws.Rows("2:3").Insert()

# The date columns (A, B) hold plain text like "2026-02-03", not real
# Excel dates -- force text format first so COM doesn't coerce the
# assigned strings into date serial numbers.
$ws.Range("A2:B3").NumberFormat = "@"
$ws.Range("A10:B12").NumberFormat = "@"

# --- New row 2: Snowflake partnership article ---
$ws.Range("A2").Value = "2026-02-03"
$ws.Range("B2").Value = "2026-02-02"
$ws.Range("C2").Value = "OpenAI"
$ws.Range("D2").Value = "Snowflake and OpenAI partner to bring frontier intelligence to enterprise data"
$ws.Range("E2").Value = "Snowflake와 OpenAI가 협력하여 기업 데이터에 최첨단 인텔리전스를 제공합니다."
$ws.Range("F2").Value = "https://openai.com/index/snowflake-partnership"

# --- New row 3: Introducing the Codex app ---
$ws.Range("A3").Value = "2026-02-03"
$ws.Range("B3").Value = "2026-02-02"
$ws.Range("C3").Value = "OpenAI"
$ws.Range("D3").Value = "Introducing the Codex app"
$ws.Range("E3").Value = "코덱스 앱을 소개합니다"
$ws.Range("F3").Value = "https://openai.com/index/introducing-the-codex-app"

# The insert shifted the old rows 8/9/10 (PVH / Powering tax donations /
# Introducing Prism) down to 10/11/12, but their relative order also
# changed: "Powering tax donations" and "Introducing Prism" moved up one
# slot each, and "PVH reimagines..." moved down to slot 12. Rewrite those
# three rows explicitly so the final order matches.
$ws.Range("A10").Value = "2026-01-28"
$ws.Range("B10").Value = "2026-01-27"
$ws.Range("C10").Value = "OpenAI"
$ws.Range("D10").Value = "Powering tax donations with AI powered personalized recommendations"
$ws.Range("E10").Value = "AI 기반 맞춤형 추천으로 세금 기부 지원"
$ws.Range("F10").Value = "https://openai.com/index/trustbank"

$ws.Range("A11").Value = "2026-01-28"
$ws.Range("B11").Value = "2026-01-27"
$ws.Range("C11").Value = "OpenAI"
$ws.Range("D11").Value = "Introducing Prism"
$ws.Range("E11").Value = "프리즘 소개"
$ws.Range("F11").Value = "https://openai.com/index/introducing-prism"

$ws.Range("A12").Value = "2026-01-28"
$ws.Range("B12").Value = "2026-01-27"
$ws.Range("C12").Value = "OpenAI"
$ws.Range("D12").Value = "PVH reimagines the future of fashion with OpenAI"
$ws.Range("E12").Value = "PVH는 OpenAI를 통해 패션의 미래를 재구상합니다."
$ws.Range("F12").Value = "https://openai.com/index/pvh-future-of-fashion"

Write-Output "OpenAI news updated"
